# Fruta / hortaliza, semanal
# Update dates (col D), volumes (col J) and prices (cols K, L, M, P)
# on the "Ciboulette" sheet to reflect the latest weekly refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44965
$ws.Range("J2").Value = 1120
$ws.Range("K2").Value = 2000
$ws.Range("L2").Value = 2500
$ws.Range("M2").Value = 2250
$ws.Range("P2").Value = 750

$ws.Range("D3").Value = 44951
$ws.Range("J3").Value = 800

$ws.Range("D4").Value = 44685
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 1500
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = 1750
$ws.Range("P4").Value = 583

$ws.Range("D5").Value = 45007
$ws.Range("J5").Value = 1160

$ws.Range("D6").Value = 44883
$ws.Range("J6").Value = 500

$ws.Range("D7").Value = 45013
$ws.Range("J7").Value = 1100

$ws.Range("D8").Value = 44971
$ws.Range("J8").Value = 1000

$ws.Range("D10").Value = 44992
$ws.Range("J10").Value = 1040

$ws.Range("D11").Value = 44978
$ws.Range("K11").Value = 1800
$ws.Range("M11").Value = 1900
$ws.Range("P11").Value = 633

$ws.Range("D12").Value = 44985

$ws.Range("D13").Value = 45035
$ws.Range("J13").Value = 1100

$ws.Range("D14").Value = 45034
$ws.Range("J14").Value = 1100

$ws.Range("D15").Value = 45041
$ws.Range("J15").Value = 1160

$ws.Range("D16").Value = 44881
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 1900
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = 1950
$ws.Range("P16").Value = 650

$ws.Range("D17").Value = 45006
$ws.Range("J17").Value = 1100

$ws.Range("D18").Value = 44953
$ws.Range("J18").Value = 1000

$ws.Range("D19").Value = 44911
$ws.Range("J19").Value = 700
$ws.Range("K19").Value = 1800
$ws.Range("L19").Value = 2000
$ws.Range("M19").Value = 1900
$ws.Range("P19").Value = 633

$ws.Range("D20").Value = 44970
$ws.Range("J20").Value = 800
$ws.Range("K20").Value = 2000
$ws.Range("L20").Value = 2500
$ws.Range("M20").Value = 2250
$ws.Range("P20").Value = 750

$ws.Range("D21").Value = 45028
$ws.Range("J21").Value = 1000

$ws.Range("D22").Value = 44827
$ws.Range("J22").Value = 1200

$ws.Range("D23").Value = 45020
$ws.Range("J23").Value = 1200

$ws.Range("D24").Value = 44848
$ws.Range("J24").Value = 1000
$ws.Range("K24").Value = 1500
$ws.Range("M24").Value = 1750
$ws.Range("P24").Value = 583

$ws.Range("D25").Value = 44999
$ws.Range("J25").Value = 1100
$ws.Range("K25").Value = 2000
$ws.Range("L25").Value = 2500
$ws.Range("M25").Value = 2250
$ws.Range("P25").Value = 750

$ws.Range("D26").Value = 44910
$ws.Range("J26").Value = 1000
